# LOM3066.xlsx update
# - Rewrites the "Objetivos" answer, adds a 3-row "Docentes responsaveis" block
#   (3 teacher names instead of 1), rewrites "Programa resumido", "Programa",
#   "Metodo", "Criterio", "Norma de recuperacao" and adds a new "Bibliografia"
#   answer row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

# ---------------------------------------------------------------------------
# 1. Make room: insert 3 blank rows at row 13 (pushes old rows 13-21 to 16-24)
# ---------------------------------------------------------------------------
$ws.Rows.Item(13).Resize(3).Insert()

# The inserted rows copy the formatting of the row above (row 12), which only
# had a styled column A. Column A should stay completely empty/unstyled on
# these new rows (rows 13-15 only have B/C content in the target layout), so
# clear that stray formatting back to "untouched".
$ws.Range("A13:A15").Clear() | Out-Null

# Give the new B13:C15 cells the same look (wrap text, top-aligned, B normal /
# C red) as the rest of the table by copying the format from row 10.
$ws.Range("B10:C10").Copy() | Out-Null
$ws.Range("B13:C15").PasteSpecial($xlPasteFormats) | Out-Null
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 2. Objetivos (row 10): replace the placeholder teacher name with the real
#    course objectives text (same text in B and C).
# ---------------------------------------------------------------------------
$objetivos = "Fornecer aos estudantes uma visão abrangente e interdisciplinar sobre materiais compósitos, além de mostrar as especificidades de cada matriz, sendo ela metálica, cerâmica ou polimérica. Ademais, deseja-se apresentar os fundamentos teóricos da mecânica de estruturas reforçadas e a partir de atividades práticas demostrar métodos de caracterização de materiais compósitos e como prepara-los."
$ws.Range("B10").Value = $objetivos
$ws.Range("C10").Value = $objetivos

# ---------------------------------------------------------------------------
# 3. Docentes responsaveis (rows 13-15): the three teachers, one per row.
# ---------------------------------------------------------------------------
$ws.Range("B13").Value = "519033 - Carlos Yujiro Shigue"
$ws.Range("C13").Value = "519033 - Carlos Yujiro Shigue"

$ws.Range("B14").Value = "1033242 - Fábio Herbst Florenzano"
$ws.Range("C14").Value = "1033242 - Fábio Herbst Florenzano"

$ws.Range("B15").Value = "1922320 - Sebastiao Ribeiro"
$ws.Range("C15").Value = "1922320 - Sebastiao Ribeiro"

# ---------------------------------------------------------------------------
# 4. Programa resumido (row 16, was row 13 before insert).
# ---------------------------------------------------------------------------
$programaResumido = "1.Introduçâo 2. Conceitos básicos sobre materiais compósitos, suas matrizes e seus processo de fabricação 3. Tipos de reforços 4. Compósitos nanoestruturados, naturais e híbridos 5. Mecânica da estrutura reforçada 6. Atividade prática"
$ws.Range("B16").Value = $programaResumido
$ws.Range("C16").Value = $programaResumido

# ---------------------------------------------------------------------------
# 5. Programa (row 18, was row 15 before insert).
# ---------------------------------------------------------------------------
$programa = "1. Conceitos básicos sobre materiais compósitos: compósitos de matriz metálica (CMM), compósitos de matriz cerâmicos (CMC) e compósitos de matriz polimérica (CMP) e nanocompósitos. 2. Tipos de Reforços: Reforços particulados, fibras curtas, fibras longas, mantas, tecidos e preformas. 3. Conceitos de Interface4. Compósitos de matriz metálica: características e processos de fabricação. 5. Compósitos de matriz cerâmica: características e processos de fabricação. 6. Compósitos de matriz polimérica: matrizes termoplásticas e termorrígidas, características físicas e químicas e processos de fabricação. 7. Compósitos nanoestruturados. 8. Compósitos Naturais. 9. Compósitos Híbridos 10. Mecânica de estruturas reforçadas. Conteúdo prático: 1. Caracterização e análise de compósitos de matriz metálica. 2. Preparação e caracterização de compósitos de matriz polimérica.(Sugestão: Considerar substituir essa parte prática pela realização do PBL descrito no item 3) 3. Visita a empresa produtora de compósitos e aulas especiais e/ou palestras com professores/pesquisadores convidados"
$ws.Range("B18").Value = $programa
$ws.Range("C18").Value = $programa

# ---------------------------------------------------------------------------
# 6. Metodo (row 21, was row 18 before insert).
# ---------------------------------------------------------------------------
$metodo = "De acordo com a atual ementa da disciplina propõe-se o uso de uma nova metodologia de ensino com o intuito de abordar o conteúdo de forma mais prática e contextualizada para que o aluno consiga relacionar os conhecimentos teóricos vistos em sala de aula com as outras disciplinas do curso. Assim, avaliação do aluno será feita através de uma prova escrita e por uma apresentação final com base nas atividades práticas desenvolvidas."
$ws.Range("B21").Value = $metodo
$ws.Range("C21").Value = $metodo

# ---------------------------------------------------------------------------
# 7. Criterio (row 22, was row 19 before insert).
# ---------------------------------------------------------------------------
$criterio = "A nota final será calculada como descrita a seguir: NF= (0,4*Avaliação escrita + 0,6 *Apresentação final)"
$ws.Range("B22").Value = $criterio
$ws.Range("C22").Value = $criterio

# ---------------------------------------------------------------------------
# 8. Norma de recuperacao (row 23, was row 20 before insert).
# ---------------------------------------------------------------------------
$norma = "Devido a cunho prático da disciplina não haverá recuperação."
$ws.Range("B23").Value = $norma
$ws.Range("C23").Value = $norma

# ---------------------------------------------------------------------------
# 9. Bibliografia (row 24, was row 21 before insert): add the new reference
#    list answer (label itself already carried over from the shift).
# ---------------------------------------------------------------------------
$bibliografia = "1. REZENDE, M. C.; COSTA, M. L.; BOTELHO, E. C. Compósitos estruturais: tecnologia e prática. São Paulo: Artliber, 2011. 396p. 2 MALLICK, P.K. Composites Engineering Handbook. New York: Marcel Dekker, 1997. 3. MATTHEWS, F.L. & RAWLINGS, R.D. Composite Materials: Engineering and Science. London: Chapman & Hall, 1994. 4. OBRAZTSOV, I.F. Mechanics of Composites. Moscow: MIR Publishers, 1982. 5. JONES R. Mechanics of Composite Materials. New York: McGraw-Hill, 1975. 6. UPADHYAYA, G.S. Sintered Metal-Ceramic Composites. Elsevier, 1984. 7. HARPER, C. A. Handbook of Plastics, Elastomers and Composites. New York: McGraw-Hill, 1992. 8. GOLDSTEIN, A.N. Handbook of Nanophase Materials. CRC Press, 1997. 9. DRESSELHAUS, M.S. Graphite Fibers and Filaments. New York: Springer-Verlag, 1988."
$ws.Range("B24").Value = $bibliografia
$ws.Range("C24").Value = $bibliografia

# ---------------------------------------------------------------------------
# 10. Column layout cleanup: column A should only be 30.71 wide on its own
#     (it used to share a (1:2) width definition with column B).
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 30.7109375

Write-Host "LOM3066 sheet rebuilt: Objetivos/Docentes/Programa/Metodo/Criterio/Norma/Bibliografia updated."
